# Auto-generated Excel COM-interop script implementing the commit:
#   "import multiline adapt basics of configuration"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

# ---- Cell values (and date number formats) for rows 1-4 ----
$ws.Range("A1").Value2 = 'ID'
$ws.Range("B1").Value2 = 'Name'
$ws.Range("C1").Value2 = 'Topic'
$ws.Range("D1").Value2 = 'Description'
$ws.Range("E1").Value2 = 'Rationale'
$ws.Range("F1").Value2 = 'Status'
$ws.Range("G1").Value2 = 'Owner'
$ws.Range("H1").Value2 = 'Invented by'
$ws.Range("I1").Value2 = 'Invented on'
$ws.Range("J1").Value2 = 'Effort estimation'
$ws.Range("K1").Value2 = 'Priority'
$ws.Range("L1").Value2 = 'Solved by'
$ws.Range("M1").Value2 = 'Type'
$ws.Range("N1").Value2 = 'Note'

$ws.Range("A2").Value2 = 'AutomaticGeneration'
$ws.Range("B2").Value2 = 'Automatic Generation of Results'
$ws.Range("C2").Value2 = 'ReqsDocument'
$ws.Range("D2").Value2 = '\textsl{rmtoo} \textbf{must} support the automatic genration of outputs.'
$ws.Range("E2").Value2 = 'Because rmtoo is aimed to be used in productive development environments, there is the need that all the different outputs (e.g. PDFs, graphs, ...) must be generated automatically (without user interaction).'
$ws.Range("F2").Value2 = 'not done'
$ws.Range("G2").Value2 = 'development'
$ws.Range("H2").Value2 = 'flonatel'
$ws.Range("I2").Value2 = 40221
$ws.Range("I2").NumberFormat = "YYYY\-MM\-DD"
$ws.Range("J2").Value2 = '3'
$ws.Range("K2").Value2 = 'development:3'
$ws.Range("L2").Value2 = 'Completed'
$ws.Range("M2").Value2 = 'requirement'

$ws.Range("A3").Value2 = 'Completed'
$ws.Range("B3").Value2 = 'Completed Requirement'
$ws.Range("C3").Value2 = 'ReqsDocument'
$ws.Range("D3").Value2 = 'It \textbf{must} be possible to check if a requirement is completed.'
$ws.Range("E3").Value2 = 'Completed means that i.e. it and all the children are finished.\par  This can be used for a ''not yet finished'' list as an output artifact. '
$ws.Range("F3").Value2 = 'finished'
$ws.Range("G3").Value2 = 'development'
$ws.Range("H3").Value2 = 'flonatel'
$ws.Range("I3").Value2 = 40243
$ws.Range("I3").NumberFormat = "YYYY\-MM\-DD"
$ws.Range("K3").Value2 = 'management:1'
$ws.Range("M3").Value2 = 'requirement'

$ws.Range("A4").Value2 = 'TestNewlines'
$ws.Range("B4").Value2 = 'Test Newlines'
$ws.Range("C4").Value2 = 'ReqsDocument'
$ws.Range("D4").Value2 = 'Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum.
ASDF'
$ws.Range("E4").Value2 = 'Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum.
QWER'
$ws.Range("F4").Value2 = 'finished'
$ws.Range("G4").Value2 = 'development'
$ws.Range("H4").Value2 = 'flonatel'
$ws.Range("K4").Value2 = 'management:1'
$ws.Range("M4").Value2 = 'requirement'
$ws.Range("N4").Value2 = 'Lipsum
Handle it well'

# ---- Wrap text for the Description/Rationale/Note columns (style index 1 in target) ----
$ws.Range("D1:D4").WrapText = $true
$ws.Range("E1:E4").WrapText = $true
$ws.Range("N1:N4").WrapText = $true

# ---- Re-apply the "ReqsDocument" style font (Arial 10) to the Topic column, including new row 4 ----
$ws.Range("C4").Font.Name = "Arial"
$ws.Range("C4").Font.Size = 10

# ---- Row heights (auto-sized by Excel originally for wrapped multi-line content) ----
$ws.Rows.Item(2).RowHeight = 91.5
$ws.Rows.Item(3).RowHeight = 57.75
$ws.Rows.Item(4).RowHeight = 215.25

# ---- Column widths (ColumnWidth = target xml width - 5/6 pixel padding) ----
$ws.Columns.Item(1).ColumnWidth = 21.39666666666667
$ws.Columns.Item(2).ColumnWidth = 21.116666666666667
$ws.Columns.Item(3).ColumnWidth = 13.606666666666666
$ws.Columns.Item(4).ColumnWidth = 35.986666666666665
$ws.Columns.Item(5).ColumnWidth = 24.866666666666667
$ws.Columns.Item(14).ColumnWidth = 22.92666666666667

# ---- Selection / scroll position ----
$ws.Range("D5").Select()

